# Insert a new weekly price record as row 89 on the "Albahaca" sheet,
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 89 (existing rows 89..176 shift to 90..177).
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new weekly record. The
# non-varying descriptive columns mirror every other row in the table.
$ws.Cells.Item(89, 1).Value = 8
$ws.Cells.Item(89, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(89, 3).Value = "Coquimbo"
$ws.Cells.Item(89, 4).Value = 45049
$ws.Cells.Item(89, 5).Value = 4
$ws.Cells.Item(89, 6).Value = 100112052
$ws.Cells.Item(89, 7).Value = "Albahaca"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 1200
$ws.Cells.Item(89, 11).Value = 2800
$ws.Cells.Item(89, 12).Value = 3000
$ws.Cells.Item(89, 13).Value = 2900
$ws.Cells.Item(89, 14).Value = "`$/paquete"
$ws.Cells.Item(89, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(89, 16).Value = 2900
$ws.Cells.Item(89, 17).Value = 1
$ws.Cells.Item(89, 18).Value = "Hortaliza"
